# scenario_settings.xlsx update:
# Adds a large block of new scenario-parameter columns (reserve, cost,
# carbon-tax and epsilon settings) to the right of the existing table on
# "Tabelle1", plus a small vertical mirror list of some of the new
# label/value pairs further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new header cells (T1:AL1) ---
$ws.Range("T1").Value = "ReserveDemand"
$ws.Range("U1").Value = "ReserveRenewables"
$ws.Range("V1").Value = "ReserveLargestUnit"
$ws.Range("W1").Value = "ReserveFast"
$ws.Range("X1").Value = "CostUnserved"
$ws.Range("Y1").Value = "CostSpilled"
$ws.Range("Z1").Value = "CostFictitiousFlow"
$ws.Range("AA1").Value = "RunningCosts"
$ws.Range("AB1").Value = "CarbonTaxy2020"
$ws.Range("AC1").Value = "CarbonTaxy2025"
$ws.Range("AD1").Value = "CarbonTaxy2030"
$ws.Range("AE1").Value = "CarbonTaxy2035"
$ws.Range("AF1").Value = "CarbonTaxy2040"
$ws.Range("AG1").Value = "CarbonTaxy2045"
$ws.Range("AH1").Value = "CarbonTaxy2050"
$ws.Range("AI1").Value = "EpsilonTransmission"
$ws.Range("AJ1").Value = "EpsilonHydropeaking"
$ws.Range("AK1").Value = "EpsilonGHG"
$ws.Range("AL1").Value = "EpsilonPMatter"

# --- Row 2: new value cells (T2:AL2) ---
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 10000
$ws.Range("Y2").Value = 0.001
$ws.Range("Z2").Value = 95000
$ws.Range("AA2").Value = 94
$ws.Range("AB2").Value = 4.5
$ws.Range("AC2").Value = 8.6
$ws.Range("AD2").Value = 12.7
$ws.Range("AE2").Value = 16.8
$ws.Range("AF2").Value = 20.9
$ws.Range("AG2").Value = 25
$ws.Range("AH2").Value = 29.1
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 0

# Cells that carry the "0.0" custom number format (same style already used
# by J2/K2 in the original workbook). Applied per contiguous block since
# multi-area (comma) ranges only honour the first area here.
$ws.Range("P2").NumberFormat = "0.0"
$ws.Range("T2:X2").NumberFormat = "0.0"
$ws.Range("Z2:AA2").NumberFormat = "0.0"
$ws.Range("AI2:AL2").NumberFormat = "0.0"

# --- Row 4: new unit / sub-header cells ---
$ws.Range("T4").Value = "factor of totsal demand"
$ws.Range("U4").Value = "factor of installed MW capacity"
$ws.Range("V4").Value = "MW"
$ws.Range("W4").Value = "MW"
$ws.Range("AB4").Value = '$ / tCO2'
$ws.Range("AJ4").Value = "Unused??"
$ws.Range("AL4").Value = "Unused??"

# --- Row 5: new annotation cells ---
$ws.Range("U5").Value = "e.g: Portugal 0.1"
$ws.Range("AB5").Value = "Source: fromChile Govt. link:http://www.minenergia.cl/mesa-geotermia/wp-content/uploads/2018/07/Informe-Final-Mesa-Geotermia.pdf"

# --- Rows 18-26: vertical label/value mirror table ---
$ws.Range("O18").Value = "CostUnserved"
$ws.Range("P18").Value = 10000
$ws.Range("O19").Value = "CostSpilled"
$ws.Range("P19").Value = 0.001
$ws.Range("O20").Value = "CostFictitiousFlows"
$ws.Range("P20").Value = 95000
$ws.Range("O21").Value = "RunningCosts"
$ws.Range("P21").Value = 94
$ws.Range("O22").Value = "CarbonTax"
$ws.Range("P22").Value = 4.5
$ws.Range("O23").Value = "EpsilonTransmission"
$ws.Range("P23").Value = 0
$ws.Range("O24").Value = "EpsilonHydropeaking"
$ws.Range("P24").Value = 0
$ws.Range("O25").Value = "EpsilonGHG"
$ws.Range("P25").Value = 0
$ws.Range("O26").Value = "EpsilonPMatter"
$ws.Range("P26").Value = 0

# --- Column widths for the newly used columns (X:AL), chosen to land as
# close as possible on the bestFit widths Excel computed for the new text. ---
$ws.Columns.Item(24).ColumnWidth = 11.5               # X
$ws.Columns.Item(25).ColumnWidth = 9                  # Y
$ws.Columns.Item(26).ColumnWidth = 16                 # Z
$ws.Columns.Item(27).ColumnWidth = 11                 # AA
$ws.Columns.Item(28).ColumnWidth = 9                  # AB
$ws.Columns.Item(29).ColumnWidth = 9                  # AC
$ws.Columns.Item(30).ColumnWidth = 9                  # AD
$ws.Columns.Item(31).ColumnWidth = 9                  # AE
$ws.Columns.Item(32).ColumnWidth = 9                  # AF
$ws.Columns.Item(33).ColumnWidth = 9                  # AG
$ws.Columns.Item(34).ColumnWidth = 9                  # AH
$ws.Columns.Item(35).ColumnWidth = 16.66666666666667  # AI
$ws.Columns.Item(36).ColumnWidth = 17.16666666666667  # AJ
$ws.Columns.Item(37).ColumnWidth = 9.666666666666666  # AK
$ws.Columns.Item(38).ColumnWidth = 12.66666666666667  # AL

# --- View state: scroll the window over and move the active selection ---
$ws.Range("Z3").Select()
